# Auto-generated edit script for Dhh-Boc.xlsx natmi update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 and add new rows 11-13 reflecting the new "M2" cluster
# and updated NATMI edge-weight statistics (per commit message: "Natmi following Dr Hou advice").

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Dhh"
$ws.Cells.Item(2,3).Value = "Boc"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 1.569424333333333
$ws.Cells.Item(2,8).Value = 4.708273
$ws.Cells.Item(2,9).Value = 0.2941318537634731
$ws.Cells.Item(2,10).Value = 0.2941318537634731
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 1.831493666666667
$ws.Cells.Item(2,14).Value = 5.494481
$ws.Cells.Item(2,15).Value = 0.04079247358478674
$ws.Cells.Item(2,16).Value = 0.04079247358478674
$ws.Cells.Item(2,17).Value = 2.874390726812556
$ws.Cells.Item(2,18).Value = 25.869516541313
$ws.Cells.Item(2,19).Value = 0.01199836587509083
$ws.Cells.Item(2,20).Value = 0.01199836587509083

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Dhh"
$ws.Cells.Item(3,3).Value = "Boc"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 1.569424333333333
$ws.Cells.Item(3,8).Value = 4.708273
$ws.Cells.Item(3,9).Value = 0.2941318537634731
$ws.Cells.Item(3,10).Value = 0.2941318537634731
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 36.96800500000001
$ws.Cells.Item(3,14).Value = 110.904015
$ws.Cells.Item(3,15).Value = 0.8233806072555885
$ws.Cells.Item(3,16).Value = 0.8233806072555884
$ws.Cells.Item(3,17).Value = 58.01848660178834
$ws.Cells.Item(3,18).Value = 522.1663794160951
$ws.Cells.Item(3,19).Value = 0.2421824643649804
$ws.Cells.Item(3,20).Value = 0.2421824643649804

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Dhh"
$ws.Cells.Item(4,3).Value = "Boc"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 1.569424333333333
$ws.Cells.Item(4,8).Value = 4.708273
$ws.Cells.Item(4,9).Value = 0.2941318537634731
$ws.Cells.Item(4,10).Value = 0.2941318537634731
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 6.098334333333334
$ws.Cells.Item(4,14).Value = 18.295003
$ws.Cells.Item(4,15).Value = 0.1358269191596248
$ws.Cells.Item(4,16).Value = 0.1358269191596247
$ws.Cells.Item(4,17).Value = 9.570874295535445
$ws.Cells.Item(4,18).Value = 86.137868659819
$ws.Cells.Item(4,19).Value = 0.03995102352340184
$ws.Cells.Item(4,20).Value = 0.03995102352340183

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Dhh"
$ws.Cells.Item(5,3).Value = "Boc"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 2.444496
$ws.Cells.Item(5,8).Value = 7.333488
$ws.Cells.Item(5,9).Value = 0.4581324022613356
$ws.Cells.Item(5,10).Value = 0.4581324022613355
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 1.831493666666667
$ws.Cells.Item(5,14).Value = 5.494481
$ws.Cells.Item(5,15).Value = 0.04079247358478674
$ws.Cells.Item(5,16).Value = 0.04079247358478674
$ws.Cells.Item(5,17).Value = 4.477078942192001
$ws.Cells.Item(5,18).Value = 40.29371047972801
$ws.Cells.Item(5,19).Value = 0.01868835391758043
$ws.Cells.Item(5,20).Value = 0.01868835391758042

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Dhh"
$ws.Cells.Item(6,3).Value = "Boc"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 2.444496
$ws.Cells.Item(6,8).Value = 7.333488
$ws.Cells.Item(6,9).Value = 0.4581324022613356
$ws.Cells.Item(6,10).Value = 0.4581324022613355
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 36.96800500000001
$ws.Cells.Item(6,14).Value = 110.904015
$ws.Cells.Item(6,15).Value = 0.8233806072555885
$ws.Cells.Item(6,16).Value = 0.8233806072555884
$ws.Cells.Item(6,17).Value = 90.36814035048002
$ws.Cells.Item(6,18).Value = 813.3132631543201
$ws.Cells.Item(6,19).Value = 0.3772173355774
$ws.Cells.Item(6,20).Value = 0.3772173355773999

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Dhh"
$ws.Cells.Item(7,3).Value = "Boc"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 2.444496
$ws.Cells.Item(7,8).Value = 7.333488
$ws.Cells.Item(7,9).Value = 0.4581324022613356
$ws.Cells.Item(7,10).Value = 0.4581324022613355
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 6.098334333333334
$ws.Cells.Item(7,14).Value = 18.295003
$ws.Cells.Item(7,15).Value = 0.1358269191596248
$ws.Cells.Item(7,16).Value = 0.1358269191596247
$ws.Cells.Item(7,17).Value = 14.907353884496
$ws.Cells.Item(7,18).Value = 134.166184960464
$ws.Cells.Item(7,19).Value = 0.06222671276635512
$ws.Cells.Item(7,20).Value = 0.0622267127663551

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Dhh"
$ws.Cells.Item(8,3).Value = "Boc"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2.0
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.165334
$ws.Cells.Item(8,8).Value = 0.4960020000000001
$ws.Cells.Item(8,9).Value = 0.03098588117774611
$ws.Cells.Item(8,10).Value = 0.03098588117774611
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 1.831493666666667
$ws.Cells.Item(8,14).Value = 5.494481
$ws.Cells.Item(8,15).Value = 0.04079247358478674
$ws.Cells.Item(8,16).Value = 0.04079247358478674
$ws.Cells.Item(8,17).Value = 0.3028081738846667
$ws.Cells.Item(8,18).Value = 2.725273564962
$ws.Cells.Item(8,19).Value = 0.001263990739444549
$ws.Cells.Item(8,20).Value = 0.001263990739444549

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Dhh"
$ws.Cells.Item(9,3).Value = "Boc"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2.0
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.165334
$ws.Cells.Item(9,8).Value = 0.4960020000000001
$ws.Cells.Item(9,9).Value = 0.03098588117774611
$ws.Cells.Item(9,10).Value = 0.03098588117774611
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 36.96800500000001
$ws.Cells.Item(9,14).Value = 110.904015
$ws.Cells.Item(9,15).Value = 0.8233806072555885
$ws.Cells.Item(9,16).Value = 0.8233806072555884
$ws.Cells.Item(9,17).Value = 6.112068138670002
$ws.Cells.Item(9,18).Value = 55.00861324803002
$ws.Cells.Item(9,19).Value = 0.0255131736604821
$ws.Cells.Item(9,20).Value = 0.02551317366048209

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Dhh"
$ws.Cells.Item(10,3).Value = "Boc"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 2.0
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.165334
$ws.Cells.Item(10,8).Value = 0.4960020000000001
$ws.Cells.Item(10,9).Value = 0.03098588117774611
$ws.Cells.Item(10,10).Value = 0.03098588117774611
$ws.Cells.Item(10,11).Value = 3.0
$ws.Cells.Item(10,12).Value = 1.0
$ws.Cells.Item(10,13).Value = 6.098334333333334
$ws.Cells.Item(10,14).Value = 18.295003
$ws.Cells.Item(10,15).Value = 0.1358269191596248
$ws.Cells.Item(10,16).Value = 0.1358269191596247
$ws.Cells.Item(10,17).Value = 1.008262008667334
$ws.Cells.Item(10,18).Value = 9.074358078006002
$ws.Cells.Item(10,19).Value = 0.004208716777819459
$ws.Cells.Item(10,20).Value = 0.004208716777819459

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Dhh"
$ws.Cells.Item(11,3).Value = "Boc"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 1.156530666666667
$ws.Cells.Item(11,8).Value = 3.469592
$ws.Cells.Item(11,9).Value = 0.2167498627974453
$ws.Cells.Item(11,10).Value = 0.2167498627974453
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 1.831493666666667
$ws.Cells.Item(11,14).Value = 5.494481
$ws.Cells.Item(11,15).Value = 0.04079247358478674
$ws.Cells.Item(11,16).Value = 0.04079247358478674
$ws.Cells.Item(11,17).Value = 2.118178591305778
$ws.Cells.Item(11,18).Value = 19.063607321752
$ws.Cells.Item(11,19).Value = 0.008841763052670938
$ws.Cells.Item(11,20).Value = 0.008841763052670938

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Dhh"
$ws.Cells.Item(12,3).Value = "Boc"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = 1.156530666666667
$ws.Cells.Item(12,8).Value = 3.469592
$ws.Cells.Item(12,9).Value = 0.2167498627974453
$ws.Cells.Item(12,10).Value = 0.2167498627974453
$ws.Cells.Item(12,11).Value = 3.0
$ws.Cells.Item(12,12).Value = 1.0
$ws.Cells.Item(12,13).Value = 36.96800500000001
$ws.Cells.Item(12,14).Value = 110.904015
$ws.Cells.Item(12,15).Value = 0.8233806072555885
$ws.Cells.Item(12,16).Value = 0.8233806072555884
$ws.Cells.Item(12,17).Value = 42.75463146798667
$ws.Cells.Item(12,18).Value = 384.7916832118801
$ws.Cells.Item(12,19).Value = 0.178467633652726
$ws.Cells.Item(12,20).Value = 0.178467633652726

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Dhh"
$ws.Cells.Item(13,3).Value = "Boc"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = 1.156530666666667
$ws.Cells.Item(13,8).Value = 3.469592
$ws.Cells.Item(13,9).Value = 0.2167498627974453
$ws.Cells.Item(13,10).Value = 0.2167498627974453
$ws.Cells.Item(13,11).Value = 3.0
$ws.Cells.Item(13,12).Value = 1.0
$ws.Cells.Item(13,13).Value = 6.098334333333334
$ws.Cells.Item(13,14).Value = 18.295003
$ws.Cells.Item(13,15).Value = 0.1358269191596248
$ws.Cells.Item(13,16).Value = 0.1358269191596247
$ws.Cells.Item(13,17).Value = 7.052910672086222
$ws.Cells.Item(13,18).Value = 63.476196048776
$ws.Cells.Item(13,19).Value = 0.02944046609204836
$ws.Cells.Item(13,20).Value = 0.02944046609204836

Write-Host "Updated sheet1 with new M2 cluster rows and refreshed NATMI statistics"
